# Auto-generated script to update transition-matrix probabilities
# reflecting additional simulated games (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2375
$ws.Cells.Item(2, 3).Value = 0.490625
$ws.Cells.Item(2, 10).Value = 0.021875
$ws.Cells.Item(2, 16).Value = 0.165625
$ws.Cells.Item(2, 19).Value = 0.08437500000000001
$ws.Cells.Item(3, 2).Value = 0.01807228915662651
$ws.Cells.Item(3, 3).Value = 0.03012048192771084
$ws.Cells.Item(3, 10).Value = 0.03012048192771084
$ws.Cells.Item(3, 16).Value = 0.7409638554216867
$ws.Cells.Item(3, 19).Value = 0.1807228915662651
$ws.Cells.Item(4, 10).Value = 0.1041666666666667
$ws.Cells.Item(4, 15).Value = 0.02083333333333333
$ws.Cells.Item(4, 16).Value = 0.6458333333333334
$ws.Cells.Item(4, 19).Value = 0.2291666666666667
$ws.Cells.Item(6, 2).Value = 0.07650273224043716
$ws.Cells.Item(6, 4).Value = 0.01092896174863388
$ws.Cells.Item(6, 6).Value = 0.06557377049180328
$ws.Cells.Item(6, 10).Value = 0.2240437158469945
$ws.Cells.Item(6, 15).Value = 0.01639344262295082
$ws.Cells.Item(6, 17).Value = 0.180327868852459
$ws.Cells.Item(6, 18).Value = 0.06010928961748634
$ws.Cells.Item(6, 19).Value = 0.366120218579235
$ws.Cells.Item(7, 2).Value = 0.1059602649006623
$ws.Cells.Item(7, 4).Value = 0.006622516556291391
$ws.Cells.Item(7, 6).Value = 0.07947019867549669
$ws.Cells.Item(7, 10).Value = 0.1390728476821192
$ws.Cells.Item(7, 15).Value = 0.03973509933774835
$ws.Cells.Item(7, 17).Value = 0.1986754966887417
$ws.Cells.Item(7, 18).Value = 0.07947019867549669
$ws.Cells.Item(7, 19).Value = 0.3509933774834437
$ws.Cells.Item(8, 2).Value = 0.0916030534351145
$ws.Cells.Item(8, 4).Value = 0.01908396946564886
$ws.Cells.Item(8, 5).Value = 0.003816793893129771
$ws.Cells.Item(8, 6).Value = 0.05725190839694656
$ws.Cells.Item(8, 10).Value = 0.133587786259542
$ws.Cells.Item(8, 15).Value = 0.02290076335877863
$ws.Cells.Item(8, 17).Value = 0.232824427480916
$ws.Cells.Item(8, 18).Value = 0.1145038167938931
$ws.Cells.Item(8, 19).Value = 0.3244274809160305
$ws.Cells.Item(9, 2).Value = 0.1633663366336634
$ws.Cells.Item(9, 4).Value = 0.009900990099009901
$ws.Cells.Item(9, 5).Value = 0.004950495049504951
$ws.Cells.Item(9, 6).Value = 0.05445544554455446
$ws.Cells.Item(9, 10).Value = 0.103960396039604
$ws.Cells.Item(9, 15).Value = 0.0396039603960396
$ws.Cells.Item(9, 17).Value = 0.1732673267326733
$ws.Cells.Item(9, 18).Value = 0.06930693069306931
$ws.Cells.Item(9, 19).Value = 0.3811881188118812
$ws.Cells.Item(10, 2).Value = 0.1200941915227629
$ws.Cells.Item(10, 4).Value = 0.03061224489795918
$ws.Cells.Item(10, 5).Value = 0.001569858712715856
$ws.Cells.Item(10, 6).Value = 0.06750392464678179
$ws.Cells.Item(10, 10).Value = 0.1122448979591837
$ws.Cells.Item(10, 15).Value = 0.01805337519623234
$ws.Cells.Item(10, 17).Value = 0.2409733124018838
$ws.Cells.Item(10, 18).Value = 0.08948194662480377
$ws.Cells.Item(10, 19).Value = 0.3194662480376766
$ws.Cells.Item(11, 7).Value = 0.1102040816326531
$ws.Cells.Item(11, 10).Value = 0.1224489795918367
$ws.Cells.Item(11, 11).Value = 0.1755102040816326
$ws.Cells.Item(11, 12).Value = 0.5836734693877551
$ws.Cells.Item(11, 19).Value = 0.00816326530612245
$ws.Cells.Item(12, 7).Value = 0.7218543046357616
$ws.Cells.Item(12, 10).Value = 0.1854304635761589
$ws.Cells.Item(12, 11).Value = 0.01324503311258278
$ws.Cells.Item(12, 12).Value = 0.05298013245033113
$ws.Cells.Item(12, 19).Value = 0.02649006622516556
$ws.Cells.Item(13, 7).Value = 0.68
$ws.Cells.Item(13, 10).Value = 0.32
$ws.Cells.Item(15, 6).Value = 0.01229508196721311
$ws.Cells.Item(15, 8).Value = 0.110655737704918
$ws.Cells.Item(15, 9).Value = 0.07377049180327869
$ws.Cells.Item(15, 10).Value = 0.3688524590163935
$ws.Cells.Item(15, 11).Value = 0.06147540983606557
$ws.Cells.Item(15, 14).Value = 0.004098360655737705
$ws.Cells.Item(15, 15).Value = 0.05737704918032787
$ws.Cells.Item(15, 19).Value = 0.3114754098360656
$ws.Cells.Item(16, 6).Value = 0.01515151515151515
$ws.Cells.Item(16, 8).Value = 0.1111111111111111
$ws.Cells.Item(16, 9).Value = 0.101010101010101
$ws.Cells.Item(16, 10).Value = 0.4646464646464646
$ws.Cells.Item(16, 11).Value = 0.1212121212121212
$ws.Cells.Item(16, 13).Value = 0.005050505050505051
$ws.Cells.Item(16, 15).Value = 0.03535353535353535
$ws.Cells.Item(16, 19).Value = 0.1464646464646465
$ws.Cells.Item(17, 6).Value = 0.006493506493506494
$ws.Cells.Item(17, 8).Value = 0.09740259740259741
$ws.Cells.Item(17, 9).Value = 0.1298701298701299
$ws.Cells.Item(17, 10).Value = 0.5173160173160173
$ws.Cells.Item(17, 11).Value = 0.06493506493506493
$ws.Cells.Item(17, 13).Value = 0.01082251082251082
$ws.Cells.Item(17, 15).Value = 0.07142857142857142
$ws.Cells.Item(17, 19).Value = 0.1017316017316017
$ws.Cells.Item(18, 6).Value = 0.01657458563535912
$ws.Cells.Item(18, 8).Value = 0.08839779005524862
$ws.Cells.Item(18, 9).Value = 0.1215469613259668
$ws.Cells.Item(18, 10).Value = 0.4033149171270718
$ws.Cells.Item(18, 11).Value = 0.1160220994475138
$ws.Cells.Item(18, 13).Value = 0.01657458563535912
$ws.Cells.Item(18, 14).Value = 0.005524861878453038
$ws.Cells.Item(18, 15).Value = 0.1104972375690608
$ws.Cells.Item(18, 19).Value = 0.1215469613259668
$ws.Cells.Item(19, 6).Value = 0.01052631578947368
$ws.Cells.Item(19, 8).Value = 0.1464114832535885
$ws.Cells.Item(19, 9).Value = 0.07751196172248803
$ws.Cells.Item(19, 10).Value = 0.4267942583732057
$ws.Cells.Item(19, 11).Value = 0.1023923444976077
$ws.Cells.Item(19, 13).Value = 0.01626794258373206
$ws.Cells.Item(19, 15).Value = 0.09473684210526316
$ws.Cells.Item(19, 19).Value = 0.1253588516746411
